# Generate Report for Handoff
#
# Moves the localization status from "In Translation" to "Ready for
# handoff" and refreshes the handoff timestamps on every sheet that
# tracks this file (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---------------------------------------------------
# Column E = zh-cn status, Column F = de-de status,
# Column G = Latest HO Xliff Generate Date.
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-11-15 17:32:07"

# --- zh-cn sheet --------------------------------------------------------
# Column C = Status, Column H = Latest Handoff Datetime.
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-11-15 17:31:53"

# --- de-de sheet --------------------------------------------------------
# Column C = Status, Column H = Latest Handoff Datetime.
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-11-15 17:32:07"

# The status text grew ("In Translation" -> "Ready for handoff"), so the
# Status columns re-size to fit the new text, matching Excel's behavior
# when a cell's contents change.
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
